$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.119.96'
$ws.Range('E2').Value = '  +4.13%  '

# Row 3
$ws.Range('D3').Value = '3.048.77'
$ws.Range('E3').Value = '  +2.01%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.52%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.93%  '

# Row 7
$ws.Range('E7').Value = '  -0.17%  '

# Row 8
$ws.Range('D8').Value = '3.042.15'
$ws.Range('E8').Value = '  +2.23%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.61%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.138'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.81%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.43'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +12.30%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.38%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000237'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.70%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.09%  '

# Row 15
$ws.Range('E15').Value = '  -0.19%  '

# Row 16
$ws.Range('D16').Value = '3.550.94'
$ws.Range('E16').Value = '  +1.99%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.75%  '

# Row 18
$ws.Range('D18').Value = '3.039.91'
$ws.Range('E18').Value = '  +1.73%  '

# Row 19
$ws.Range('D19').Value = '61.075.03'
$ws.Range('E19').Value = '  +3.96%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.13%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.20%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.732'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.80%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.52%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.47%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.73%  '

# Row 26
$ws.Range('E26').Value = '  +0.18%  '

# Row 27
$ws.Range('E27').Value = '  +8.08%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.53%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.99%  '

# Row 31
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.65%  '

# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.04%  '

# Row 33
$ws.Range('E33').Value = '  +7.95%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0811'
$ws.Range('E34').Value = '  +8.56%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.60%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.76%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.80%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.10%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.46%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '408.33'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.97%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0363'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.59%  '

# Row 43
$ws.Range('D43').Value = '2.780.60'
$ws.Range('E43').Value = '  +2.34%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.107'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.57%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.261'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.28%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '37.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +17.47%  '

# Row 47
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.29%  '

# Row 48
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.01%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.85%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.111'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.40%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.79%  '
